# Updated cryptos list values (price + 1h volume change) to match the refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a cell value as literal text even when it looks like a number
# (e.g. "610.00", "1.00") -- prefixing with an apostrophe mirrors typing the
# value into Excel with the cell pre-formatted/entered as Text, so values such as
# trailing zeros or exact decimal digits are preserved instead of being normalized
# into a floating point number.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
}

$ws.Range("D2").Value = '73.388.94'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '3.977.70'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue $ws.Range("D5") '610.00'
$ws.Range("E5").Value = '  +6.43%  '
Set-TextValue $ws.Range("D6") '168.69'
$ws.Range("E6").Value = '  +10.80%  '
Set-TextValue $ws.Range("D7") '0.683'
$ws.Range("E7").Value = '  -1.97%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-TextValue $ws.Range("D9") '0.790'
$ws.Range("E9").Value = '  +2.49%  '
Set-TextValue $ws.Range("D10") '0.187'
$ws.Range("E10").Value = '  +8.09%  '
Set-TextValue $ws.Range("D11") '56.26'
$ws.Range("E11").Value = '  +4.10%  '
$ws.Range("E12").Value = '  +2.32%  '
Set-TextValue $ws.Range("D13") '11.30'
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("D14").Value = '4.614.79'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '3.977.15'
$ws.Range("E15").Value = '  -2.11%  '
Set-TextValue $ws.Range("D16") '14.29'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("E17").Value = '  +1.52%  '
Set-TextValue $ws.Range("D18") '20.77'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").Value = '73.304.00'
$ws.Range("E19").Value = '  -0.03%  '
Set-TextValue $ws.Range("D20") '0.130'
$ws.Range("E20").Value = '  -1.30%  '
Set-TextValue $ws.Range("D21") '454.62'
$ws.Range("E21").Value = '  +1.75%  '
Set-TextValue $ws.Range("D22") '4.83'
$ws.Range("E22").Value = '  +6.68%  '
Set-TextValue $ws.Range("D23") '96.29'
$ws.Range("E23").Value = '  -2.03%  '
Set-TextValue $ws.Range("D24") '3.42'
$ws.Range("E24").Value = '  -4.02%  '
Set-TextValue $ws.Range("D25") '14.25'
$ws.Range("E25").Value = '  -3.16%  '
Set-TextValue $ws.Range("D26") '4.19'
$ws.Range("E26").Value = '  -2.29%  '
Set-TextValue $ws.Range("D27") '11.08'
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("E28").Value = '  +0.24%  '
Set-TextValue $ws.Range("D29") '10.55'
$ws.Range("E29").Value = '  -4.61%  '
Set-TextValue $ws.Range("D30") '36.40'
$ws.Range("E30").Value = '  -2.37%  '
Set-TextValue $ws.Range("D31") '7.96'
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("E33").Value = '  +16.28%  '
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("E35").Value = '  -1.03%  '
Set-TextValue $ws.Range("D36") '70.95'
$ws.Range("E36").Value = '  +3.83%  '
Set-TextValue $ws.Range("D37") '648.53'
$ws.Range("E37").Value = '  -5.68%  '
Set-TextValue $ws.Range("D38") '0.432'
$ws.Range("E38").Value = '  -3.68%  '
$ws.Range("E39").Value = '  +0.14%  '
Set-TextValue $ws.Range("D40") '0.147'
$ws.Range("E40").Value = '  -1.67%  '
Set-TextValue $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("E42").Value = '  +0.12%  '
Set-TextValue $ws.Range("D43") '3.25'
$ws.Range("E43").Value = '  +41.21%  '
Set-TextValue $ws.Range("D44") '0.0484'
$ws.Range("E44").Value = '  -2.63%  '
Set-TextValue $ws.Range("D45") '10.63'
$ws.Range("E45").Value = '  -5.41%  '
Set-TextValue $ws.Range("D46") '3.13'
$ws.Range("E46").Value = '  -6.12%  '
$ws.Range("E48").Value = '  +8.62%  '
Set-TextValue $ws.Range("D49") '3.47'
$ws.Range("E49").Value = '  +3.16%  '
Set-TextValue $ws.Range("D50") '2.59'
$ws.Range("E50").Value = '  -4.72%  '
$ws.Range("E51").Value = '  -3.13%  '
